# Weekly update: insert a new price record as row 398, pushing the
# existing rows 398-423 down to 399-424 (dimension grows to A1:R424).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 398 - this shifts rows
# 398..423 down to 399..424 and keeps their content/formatting intact.
$ws.Rows.Item(398).Insert()

# Populate the new row 398 with this week's record. Columns that stay
# the same as the (now shifted-down) neighbouring rows are carried over
# explicitly so the row is fully formed.
$ws.Cells.Item(398, 1).Value2  = 3                            # A: Mercado ID
$ws.Cells.Item(398, 2).Value2  = "Femacal de La Calera"       # B: Mercado
$ws.Cells.Item(398, 3).Value2  = "Coquimbo"                   # C: Región
$ws.Cells.Item(398, 4).Value2  = 44826                        # D: Fecha
$ws.Cells.Item(398, 5).Value2  = 5                            # E: Codreg
$ws.Cells.Item(398, 6).Value2  = 100112043                    # F: Categoría ID
$ws.Cells.Item(398, 7).Value2  = "Pepino ensalada"             # G: Categoría
$ws.Cells.Item(398, 8).Value2  = "Sin especificar"             # H: Variedad
$ws.Cells.Item(398, 9).Value2  = "Primera"                     # I: Calidad
$ws.Cells.Item(398, 10).Value2 = 125                           # J: Volumen
$ws.Cells.Item(398, 11).Value2 = 17000                         # K: Precio mínimo
$ws.Cells.Item(398, 12).Value2 = 18000                         # L: Precio máximo
$ws.Cells.Item(398, 13).Value2 = 17520                         # M: Precio promedio ponderado
$ws.Cells.Item(398, 14).Value2 = "$/caja 60 unidades"          # N: Unidad de comercialización
$ws.Cells.Item(398, 15).Value2 = "Región de Arica y Parinacota" # O: Origen
$ws.Cells.Item(398, 16).Value2 = 292                           # P: Precio $/Kg
$ws.Cells.Item(398, 17).Value2 = 60                            # Q: Kg o Unidades
$ws.Cells.Item(398, 18).Value2 = "Hortaliza"                   # R: Clasificación

# Match the date-formatted number format already used by column D in
# this table (carried by the row that was pushed down from 398 to 399).
$ws.Cells.Item(398, 4).NumberFormat = $ws.Cells.Item(399, 4).NumberFormat
